# Applies the "fixed typos in manuscript and response" edit.
$d = $word.ActiveDocument

# 1. Update the date field's cached result text (literal date string).
$d.Content.Find.Execute("Wednesday, 20 March 2019", $false, $false, $false, $false, $false,
                         $true, 1, $false, "Friday, 22 March 2019", 2) | Out-Null

# 2. Merge the split "22 October 2018" runs back into a single run
#    (visible text is unchanged, just the run-splitting typo is fixed).
$d.Content.Find.Execute("22 October 2018", $false, $false, $false, $false, $false,
                         $true, 1, $false, "22 October 2018", 2) | Out-Null

# 3. Merge the split "reviewer's concerns" runs back into a single run.
$d.Content.Find.Execute("reviewer" + [char]0x2019 + "s concerns", $false, $false, $false, $false, $false,
                         $true, 1, $false, "reviewer" + [char]0x2019 + "s concerns", 2) | Out-Null

# 4. Add the second author after "Thomas G. Close, Ph. D."
$d.Content.Find.Execute("Thomas G. Close, Ph. D.", $false, $false, $false, $false, $false,
                         $true, 1, $false, "Thomas G. Close, Ph. D. and Gary F. Egan, Ph. D.", 2) | Out-Null
